$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("A995")
$r.Value = "TESTROW"
$r.Font.Bold = $false
$r.Font.Size = 7
$r.HorizontalAlignment = -4108  # xlCenter
$r.VerticalAlignment = -4108
$r.WrapText = $false
$r.Borders.LineStyle = 1
$r.Borders.ColorIndex = 64
Write-Host "done"
